$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "updated date" cell
$ws.Range("A1").Value = "更新日期：2025.01.31 13:17:43"

# 2. Fix the confidence marker on D281
$ws.Range("D281").Value = "*maa://45842"

# 3. Delete the duplicate "涤火杰西卡" row (row 310); everything below shifts up
$ws.Rows(310).Delete()

# 4. Update 玛露西尔's row (now row 351) D/C values to add a second maa link
$ws.Range("C351").Value = "2"
$ws.Range("D351").Value = "maa://41110, maa://45605"
